# Apply FX data update to Guatemala_FX sheet:
# - Fix row 312 (D/E/F values)
# - Append rows 313-315 with new monthly FX data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct existing row 312 values (high/low/close)
$ws.Range("D312").Value = 7.8715
$ws.Range("E312").Value = 7.8045
$ws.Range("F312").Value = 7.8712

# New row 313
$ws.Range("A313").Value = 45170.33333333334
$ws.Range("B313").Value = "FX_IDC:USDGTQ"
$ws.Range("C313").Value = 7.8712
$ws.Range("D313").Value = 7.874
$ws.Range("E313").Value = 7.858
$ws.Range("F313").Value = 7.858
$ws.Range("G313").Value = 0

# New row 314
$ws.Range("A314").Value = 45201.375
$ws.Range("B314").Value = "FX_IDC:USDGTQ"
$ws.Range("C314").Value = 7.858
$ws.Range("D314").Value = 7.858
$ws.Range("E314").Value = 7.803
$ws.Range("F314").Value = 7.8325
$ws.Range("G314").Value = 0

# New row 315
$ws.Range("A315").Value = 45231.375
$ws.Range("B315").Value = "FX_IDC:USDGTQ"
$ws.Range("C315").Value = 7.8325
$ws.Range("D315").Value = 7.8325
$ws.Range("E315").Value = 7.824
$ws.Range("F315").Value = 7.8245
$ws.Range("G315").Value = 0

# Match formatting of existing date column (style copied from A312)
$ws.Range("A312").Copy() | Out-Null
$ws.Range("A313:A315").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
